$d = $word.ActiveDocument

# --- Paragraph 1: intro paragraph -----------------------------------------
$d.Content.Find.Execute(
  "El diseño de nivel de “Marshallow: Pilferage in Yolk-Town” es algo que se ha tenido más o menos claro desde un primer momento, ya que se quería crear un escenario no demasiado grande para concentrar la acción de robo del ladrón y a la vez un escenario lo suficientemente amplio como para dar sensación de libertad al jugador, pudiendo este moverse libremente.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "El diseño de nivel de “Marshallow: Pilferage in YolkTown” es algo que se ha tenido claro desde prácticamente el primer momento, ya que se quería crear un escenario no demasiado grande para concentrar la acción de robo del ladrón y a la vez un lo suficientemente amplio como para dar sensación de libertad al jugador.",
  2) | Out-Null

# --- Paragraph 2: zones / distribution paragraph ---------------------------
$d.Content.Find.Execute(
  "En cuanto al escenario, se ha decidido diseñar un nivel en el que encontramos 6 zonas de interés, distribuidas una de ellas en el centro del escenario (la plaza del pueblo) y las cinco restantes alrededor de esta, creando así una jugabilidad más frenética al no haber mucha distancia entre una zona y otra, intentando así no alargar demasiado las partidas (sucedería si el escenario fuese mucho más grande o distribuido más ampliamente en vez de todo tan concentrado). ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Se ha decidido diseñar un nivel en el que encontramos 6 zonas de interés distribuidas en forma de estrella con una de ellas en el centro (la plaza del pueblo) y las cinco restantes alrededor de esta. Al ser un nivel compacto, se genera una jugabilidad más frenética ya que no hay mucha distancia entre una zona y otra, intentando así no alargar demasiado las partidas. ",
  2) | Out-Null

# --- Paragraph 3: aldeanos distribution paragraph ---------------------------
$d.Content.Find.Execute(
  "Esta distribución se diseña de tal forma para que haya suficientes aldeanos por zona y también se implementará que no pueda haber más de X aldeanos por zona permitiendo así que no haya zonas vacías en el escenario y no aburrir al jugador. Además, el ladrón tratará de robar en las zonas de interés para que así el jugador sepa donde investigar. ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Esta distribución se aprovecha de tal forma que el número de aldeanos por zona sea equitativo, es decir, que no se dé el caso de que haya un exceso de habitantes o falta de ellos en una zona. También se tendrá en cuenta la distribución a la hora de que el ladrón decida robar, ya que por norma general lo hará dentro de dichas zonas (esto se explicará más a fondo en la descripción detallada de los agentes). ",
  2) | Out-Null

# --- Paragraph 4: placeholder/image description paragraph ------------------
$d.Content.Find.Execute(
  "A continuación, se muestra un placeholder del escenario usado para pruebas, resaltando en rojo las zonas comentadas anteriormente y pudiendo observar como encontramos múltiples conexiones entre ellas.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "A continuación, se muestra una imagen en vista cenital del blocking del escenario, resaltando en rojo las zonas comentadas anteriormente. Se puede observar cómo encontramos múltiples conexiones entre estas, generando así un gameplay fluido.",
  2) | Out-Null

# --- Last paragraph: entradas/salidas paragraph -----------------------------
$d.Content.Find.Execute(
  "Por último, comentar que hemos decidido crear más de una entrada/salida en cada zona para así dificultar un poco la búsqueda del ladrón ya que, si solo hubiese una, al ver a alguien huir habiendo una víctima en dicha zona sabremos que ese pueblerino huyendo será el culpable del robo.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Cabe comentar que hemos interconectado todas las zonas desde varios puntos para así aumentar la dificultad a la hora de encontrar al ladrón, debido a que de darse el caso en el que sólo hay un pasillo de salida desde una de las zonas, el ladrón solo puede haber escapado por este y cabe la posibilidad de que se encuentre al jugador, haciendo muy fácil y aburrida la experiencia de juego.",
  2) | Out-Null

# That last paragraph also gained paragraph-level justification (w:jc="both"),
# matching every other body paragraph in the document.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.ParagraphFormat.Alignment = 3
